$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 changes from text "0948579427" to numeric value 948579427
$ws.Range("D3").Value = 948579427

# Add new row 4 data
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "nguyenminh"
$ws.Range("C4").Value = "nguynmin3@gmail.com"
$ws.Range("D4").Value = "'123456789"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "Formula"
$ws.Range("F4").Value = "scrypt:32768:8:1`$cWWwvrIaJzYDQl0M`$c5fb914ec2ecb440e623f24fda44dae20d4e15e51d9a883bb9ca9e5f30314cf715cac62eb6e70df7788350c390500c38c8cb017c192d070953c788c9d7a10bb5"
